$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# The shelter entry "Abangan Sur Covered Court" (row 7, all-zero Lat/Long) was a
# duplicate/bad entry — remove the whole row, shifting everything below it up.
$ws.Rows.Item(7).Delete()

# Rename the Lat/Long headers from the old "xDegrees"/"yDegrees" placeholders
# to proper "Latitude"/"Longitude" labels.
$ws.Range("B1").Value = "Latitude"
$ws.Range("C1").Value = "Longitude"

# Restore the active selection to A7, matching where the author left off.
$null = $ws.Range("A7").Select()

Write-Output "done"
